# MCA_Belaege_Szenario1: split the raw sheet into a data sheet ("Gewichtung")
# and a narrative sheet ("Begründung") that documents the weighting rationale.

$wb = $excel.ActiveWorkbook

# --- Rename the original (only) sheet: it now holds just the weighting table ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Gewichtung"

# --- Insert the new "Begründung" sheet right after "Gewichtung" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Begründung"

# --- Write the scenario-1 rationale paragraph into A1 ---
$text = "Szenario 1: `"Ausgewogene Nachhaltigkeit`" (Baseline)`nSzenario 1, `"Ausgewogene Nachhaltigkeit`", dient als Baseline-Szenario. Die Gewichtung (Umweltbelastung: 30%; Langlebigkeit & Wirtschaftlichkeit: 25%; Multifunktionale Nutzungsqualität: 25%; Kreislauffähigkeit: 20%) repräsentiert einen holistischen Ansatz, bei dem alle vier zentralen Bewertungsbereiche annähernd gleich stark priorisiert werden.`nDiese Verteilung legitimiert sich direkt aus Forschungsziel 1, welches eine umfassende Analyse von Umweltverträglichkeit, Lebensdauer und Kreislauffähigkeit fordert, sowie Forschungsziel 2, das die Anwendung eines breiten Bewertungsrasters impliziert. Das Szenario bildet die Referenz, an der die Sensitivität der Prioritätenverschiebung in den Folgeszenarien gemessen wird."
$ws2.Range("A1").Value = $text
$ws2.Range("A1").WrapText = $true
$ws2.Range("A1").VerticalAlignment = -4160

# --- Size the column/row so the wrapped paragraph is fully visible ---
$ws2.Columns.Item(1).ColumnWidth = 52.33
$ws2.Rows.Item(1).RowHeight = 246.5
$ws2.PageSetup.TopMargin = 56.69291338582677
$ws2.PageSetup.BottomMargin = 56.69291338582677

# --- Restore selections per sheet, and leave "Begründung" as the active tab ---
[void]$ws1.Range("B28").Select()
[void]$ws2.Range("C9").Select()
[void]$ws2.Activate()
